$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2412"
$ws.Range("F16").Value = 52000

$ws.Range("C17").Value = "1001901797"
$ws.Range("D17").Value = "ROSA GISELA MORALES MEZA"
$ws.Range("E17").Value = "2412"
$ws.Range("F17").Value = 52000

$ws.Range("C18").Value = "1001896770"
$ws.Range("D18").Value = "YURIS MILENA MORENO MEZA"

$ws.Range("C19").Value = "1006291866"
$ws.Range("D19").Value = "JUAN PABLO GUTIERREZ SALAZAR"
$ws.Range("G19").Value = 1423500

$ws.Range("G20").Value = 1423500

$ws.Range("C21").Value = "45531768"
$ws.Range("D21").Value = "LISBETH PACHECO VALENCIA"
$ws.Range("E21").Value = "2504"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("E22").Value = "2504"

$ws.Range("C23").Value = "45531768"
$ws.Range("D23").Value = "LISBETH PACHECO VALENCIA"
$ws.Range("E23").Value = "2505"

$ws.Range("E24").Value = "2505"
$ws.Range("F24").Value = 56940
